$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 1836.25
$ws.Range("I20").Value = 1836.25
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1836.25
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1606.25

# Row 35
$ws.Range("H35").Value = 1836.25
$ws.Range("I35").Value = 1836.25
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1836.25
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1457.25

# Row 137
$ws.Range("H137").Value = 1965.8511
$ws.Range("I137").Value = 3304.3076
$ws.Range("J137").Value = 1454.0883
$ws.Range("K137").Value = 9912.9228
$ws.Range("L137").Value = 4362.2649
$ws.Range("M137").Value = -7362.9228

# Row 138
$ws.Range("H138").Value = 2629.1194
$ws.Range("I138").Value = 1534.75
$ws.Range("J138").Value = 3629.6858
$ws.Range("K138").Value = 4604.25
$ws.Range("L138").Value = 10889.0574
$ws.Range("M138").Value = 535.75
$ws.Range("N138").Value = -21169.0574


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1899.9111
$ws.Range("I61").Value = 1873.2174
$ws.Range("J61").Value = 1927.8182
$ws.Range("K61").Value = 1873.2174
$ws.Range("L61").Value = 1927.8182
$ws.Range("M61").Value = -1661.2174
$ws.Range("N61").Value = -2351.8182

# Row 74
$ws.Range("H74").Value = 2682.6667
$ws.Range("I74").Value = 3297.3333
$ws.Range("J74").Value = 2477.7778
$ws.Range("K74").Value = 3297.3333
$ws.Range("L74").Value = 2477.7778
$ws.Range("M74").Value = -2423.3333
$ws.Range("N74").Value = -4225.7778

# Row 77
$ws.Range("H77").Value = 2682.6667
$ws.Range("I77").Value = 3297.3333
$ws.Range("J77").Value = 2477.7778
$ws.Range("K77").Value = 16486.6665
$ws.Range("L77").Value = 12388.889
$ws.Range("M77").Value = -12118.6665
$ws.Range("N77").Value = -21124.889

# Row 136
$ws.Range("H136").Value = 1899.9111
$ws.Range("I136").Value = 1873.2174
$ws.Range("J136").Value = 1927.8182
$ws.Range("K136").Value = 5619.6522
$ws.Range("L136").Value = 5783.4546
$ws.Range("M136").Value = -3069.6522
$ws.Range("N136").Value = -10883.4546


$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2767.7112
$ws.Range("I134").Value = 1929.7916
$ws.Range("J134").Value = 3725.3333
$ws.Range("K134").Value = 5789.3748
$ws.Range("L134").Value = 11175.9999
$ws.Range("M134").Value = -3254.3748

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 956.4167
$ws.Range("I19").Value = 537.8
$ws.Range("J19").Value = 3049.5
$ws.Range("K19").Value = 537.8
$ws.Range("L19").Value = 3049.5
$ws.Range("M19").Value = -367.8

# Row 24
$ws.Range("H24").Value = 956.4167
$ws.Range("I24").Value = 537.8
$ws.Range("J24").Value = 3049.5
$ws.Range("K24").Value = 537.8
$ws.Range("L24").Value = 3049.5
$ws.Range("M24").Value = -367.8

# Row 31
$ws.Range("H31").Value = 2433.58
$ws.Range("I31").Value = 1155.151
$ws.Range("J31").Value = 3875.2126
$ws.Range("K31").Value = 1155.151
$ws.Range("L31").Value = 3875.2126
$ws.Range("M31").Value = -860.1510000000001
$ws.Range("N31").Value = -4465.2126

# Row 33
$ws.Range("H33").Value = 5989.5
$ws.Range("I33").Value = 5989.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 5989.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -5610.5

# Row 34
$ws.Range("H34").Value = 2433.58
$ws.Range("I34").Value = 1155.151
$ws.Range("J34").Value = 3875.2126
$ws.Range("K34").Value = 1155.151
$ws.Range("L34").Value = 3875.2126
$ws.Range("M34").Value = -953.1510000000001
$ws.Range("N34").Value = -4279.2126

# Row 36
$ws.Range("H36").Value = 3994.5
$ws.Range("I36").Value = 2989
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 2989
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -2601

# Row 40
$ws.Range("H40").Value = 3994.5
$ws.Range("I40").Value = 2989
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2989
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2829

# Row 62
$ws.Range("H62").Value = 7797.5
$ws.Range("I62").Value = 2065.7693
$ws.Range("J62").Value = 22700
$ws.Range("K62").Value = 2065.7693
$ws.Range("L62").Value = 22700
$ws.Range("M62").Value = -1441.7693
$ws.Range("N62").Value = -23948

# Row 65
$ws.Range("H65").Value = 7797.5
$ws.Range("I65").Value = 2065.7693
$ws.Range("J65").Value = 22700
$ws.Range("K65").Value = 10328.8465
$ws.Range("L65").Value = 113500
$ws.Range("M65").Value = -7208.8465
$ws.Range("N65").Value = -119740

# Row 92
$ws.Range("H92").Value = 16940.4
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 16940.4
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 16940.4
$ws.Range("N92").Value = -21932.4

# Row 134
$ws.Range("H134").Value = 1483.122
$ws.Range("I134").Value = 1047.6
$ws.Range("J134").Value = 2163.625
$ws.Range("K134").Value = 3142.8
$ws.Range("L134").Value = 6490.875
$ws.Range("M134").Value = -607.7999999999997
$ws.Range("N134").Value = -11560.875


$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 104
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 104
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 312
$ws.Range("N12").Value = -658

# Row 113
$ws.Range("H113").Value = 696.0714
$ws.Range("I113").Value = 642.4167
$ws.Range("J113").Value = 736.3125
$ws.Range("K113").Value = 1927.2501
$ws.Range("L113").Value = 2208.9375
$ws.Range("M113").Value = 242.7499
$ws.Range("N113").Value = -6548.9375


$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 1526.625
$ws.Range("I9").Value = 700.8333
$ws.Range("J9").Value = 4004
$ws.Range("K9").Value = 700.8333
$ws.Range("L9").Value = 4004
$ws.Range("M9").Value = -530.8333

# Row 13
$ws.Range("H13").Value = 279
$ws.Range("I13").Value = 268.33334
$ws.Range("J13").Value = 295
$ws.Range("K13").Value = 268.33334
$ws.Range("L13").Value = 295
$ws.Range("M13").Value = -129.33334
$ws.Range("N13").Value = -573

# Row 70
$ws.Range("H70").Value = 7016
$ws.Range("I70").Value = 7291.636
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 7291.636
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -7021.636

# Row 73
$ws.Range("H73").Value = 7016
$ws.Range("I73").Value = 7291.636
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 7291.636
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -6355.636


$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 380.7143
$ws.Range("I9").Value = 294.16666
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 294.16666
$ws.Range("L9").Value = 900
$ws.Range("M9").Value = -70.16665999999998

# Row 22
$ws.Range("H22").Value = 515
$ws.Range("I22").Value = 481.66666
$ws.Range("J22").Value = 565
$ws.Range("K22").Value = 481.66666
$ws.Range("L22").Value = 565
$ws.Range("M22").Value = -186.66666
$ws.Range("N22").Value = -1155

# Row 27
$ws.Range("H27").Value = 515
$ws.Range("I27").Value = 481.66666
$ws.Range("J27").Value = 565
$ws.Range("K27").Value = 481.66666
$ws.Range("L27").Value = 565
$ws.Range("M27").Value = -374.66666
$ws.Range("N27").Value = -779

# Row 136
$ws.Range("H136").Value = 2249.65
$ws.Range("I136").Value = 2037.375
$ws.Range("J136").Value = 3098.75
$ws.Range("K136").Value = 6112.125
$ws.Range("L136").Value = 9296.25
$ws.Range("M136").Value = -3562.125
$ws.Range("N136").Value = -14396.25


$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 40143
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 40143
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 40143
$ws.Range("N123").Value = -49943

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

